$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3 (SamplesTab row): trim the "Sample ID" query down to 4 columns,
#     dropping the Tumor / Analyte Type columns ("Added CDS All studies testcase") ---
$newSamplesQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001437' AND gi.library_layout = 'Paired-End'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newSamplesQuery

# --- D3:E3 / D4:E4 no longer carry the Tsv/Web file-name cells ---
$ws.Range("D3:E3").ClearContents()
$ws.Range("D4:E4").ClearContents()

# --- Selection moves from C4 to the D3:E4 block (active cell D4) ---
$ws.Range("D3:E4").Select() | Out-Null
